$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string used in D2 ("bank bni" -> "Transfer Bank")
$ws.Range("D2").Value = "Transfer Bank"

# Update the existing row 2 account numbers
$ws.Range("C2").Value = 1234567890
$ws.Range("F2").Value = 1234567890

# Add new row 3 data (duplicate of row 2 pattern with new amounts)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2023
$ws.Range("C3").Value = 1234567890
$ws.Range("D3").Value = "Transfer Bank"
$ws.Range("E3").Value = "bank bri"
$ws.Range("F3").Value = 1234567890
$ws.Range("G3").Value = "admin"
$ws.Range("H3").Value = "idr"
$ws.Range("I3").Value = 20000000
$ws.Range("J3").Value = 20000000
